$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rapport Hors pont")

# A2 ("03250002") looks numeric, so force text formatting first or Excel
# will silently coerce it to a Number and drop the leading zero. Revert the
# cell style back to Normal afterwards so no stray number format sticks.
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "03250002"
$ws.Range("A2").Style = "Normal"

$ws.Range("C2").Value = "Cause incident 1"
$ws.Range("D2").Value = "SALLE DE CONFERENCE ROOM"
$ws.Range("E2").Value = "2NT"
$ws.Range("H2").Value = "plomb 1"
$ws.Range("J2").Value = ":ETS DJOUBISSIE ET FILS"
$ws.Range("K2").Value = "(CETOSTEARYL ALCOHOL"
$ws.Range("L2").Value = "13 BSMAT GAP MCD MOULINS"
$ws.Range("M2").Value = "O45FFG"
$ws.Range("O2").Value = "ALI"
$ws.Range("P2").Value = "455RG"
$ws.Range("Q2").Value = "Admin User"
